# Fix Excel export to exclude header in Save_Settings
# Settings_Custom.xlsx should only contain the two "setting name / value"
# rows (Rounded, Time_Delay) with no header row and no special styling
# (the old bold/centered/bordered header row is gone entirely).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop any formatting (bold font, thin border, centered alignment) that was
# applied to the old header row before we remove/overwrite the data.
$ws.UsedRange.ClearFormats()

# Remove the old rows 3-8 (the numeric header/config rows that are no
# longer exported) so only two data rows remain.
$ws.Rows("3:8").Delete()

# Rewrite the remaining two rows with the new (header-less) values.
$ws.Range("A1").Value = "Rounded"
$ws.Range("B1").Value = 2
$ws.Range("A2").Value = "Time_Delay"
$ws.Range("B2").Value = 0.5
